$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 4 ("Siswa 3") is being dropped entirely; the old row 5
# ("Ilham Shiddiq") takes its place as the new row 4. Deleting row 4
# shifts row 5 up and reproduces exactly that outcome (including K4=10,
# N4="Secret", O4="12 RPL A" which we clear below) while also updating
# the sheet dimension to A1:O4 automatically.
$ws.Rows(4).Delete()

# Row 2: "Siswa 1" -> "Hariz Sufyan Munawar"
$ws.Range("A2").Value = "Hariz Sufyan Munawar"
$ws.Range("B2").Value = "'0040172372"
$ws.Range("C2").Value = "'123213"
$ws.Range("D2").Value = "munawarhariz@gmail.com"
$ws.Range("E2").Value = "Bandung"
$ws.Range("F2").Value = "'2021-03-08"
$ws.Range("G2").Value = "Pria"
$ws.Range("I2").Value = "Lembah Teratai Blok N no.12"
$ws.Range("J2").Value = "SMP Negeri 3 Cimahi"
$ws.Range("L2").Value = "Moch. Yusuf"
$ws.Range("M2").Value = "Zulaekah"

# Row 3: "Siswa 2" -> "Shaddam Amru"
$ws.Range("A3").Value = "Shaddam Amru"
$ws.Range("B3").Value = "'0031068496"
$ws.Range("C3").Value = "'123213512"
$ws.Range("D3").Value = "shaddam.a.h@gmail.com"
$ws.Range("F3").Value = "'2021-03-10"
$ws.Range("I3").Value = "Cijerah"
$ws.Range("J3").Value = "SMP Negeri 4 Bandung"
$ws.Range("L3").Value = "Amri Hasibuan"
$ws.Range("M3").Value = "Ani Hasibuan"
$ws.Range("O3").Value = "'"

# Row 4 (was old row 5, "Ilham Shiddiq"): only the class column is cleared
$ws.Range("O4").Value = "'"
